$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power BusInfo")

# Select row 7 (Node_1 row) as the last remaining node, matching the
# sheet's recorded selection after the edit.
$ws.Range("A7:XFD7").Select()

# Delete the entire rows for Node_2 (row 8) and Node_3 (row 9), leaving
# only Node_1 and turning the sheet into a single-node, no-network setup.
# This shifts rows 10-13 up by two.
$ws.Range("A8:A9").EntireRow.Delete()

# The defined ranges used for the table / autofilter need to be shrunk
# to match the now-smaller data region.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Power BusInfo!_FilterDatabase") {
        $n.RefersTo = "='Power BusInfo'!`$M`$8:`$R`$11"
    }
    if ($n.Name -eq "businfo") {
        $n.RefersTo = "='Power BusInfo'!`$B`$3:`$K`$8"
    }
}
